# Remove the "url.cell_type" column (column C) from the "survey" sheet.
# This was used to force the adjacent "url" cell to be parsed as a formula;
# the converter now does its own expression sanity-checking, so the helper
# column is no longer needed. Deleting the column shifts every later
# column left by one and Excel automatically drops the now-unused
# "url.cell_type" / "formula" shared-string entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
$ws.Columns.Item(3).Delete()

# The "survey" sheet tab is now the active/selected one (it used to be
# "initial"), with cell B6 selected.
$ws.Activate()
[void]$ws.Range("B6").Select()
